$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$aPattern = @(10002,10003,10004,10005,10006,10007,10008,10009,10010)

for ($row = 102; $row -le 146; $row++) {
    $idx = ($row - 102) % 9
    $aVal = $aPattern[$idx]
    $bVal = 3000121 + ($row - 102)
    $ws.Cells.Item($row, 1).Value = $aVal
    $ws.Cells.Item($row, 2).Value = $bVal
    $ws.Cells.Item($row, 3).Value = "eng"
    $ws.Cells.Item($row, 4).Value = $true
    $ws.Cells.Item($row, 5).Value = "superadmin"
    $ws.Cells.Item($row, 6).Value = "now()"
}

$win = $excel.Windows.Item(1)
$win.ScrollRow = 128
$win.ScrollColumn = 1

$ws.Range("A102:F146").Select()

$ps = $ws.PageSetup
$ps.Orientation = 1

Write-Host "Done adding rows 102-146"